$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 380-428 with 2024/2025 election data (Sachsen, Thueringen,
# Brandenburg, Hamburg) - values written in the exact order needed so that
# newly introduced strings land in the shared string table in the same
# sequence as the authoritative export.
$ws.Cells.Item(380, 1).Value = 'SN'
$ws.Cells.Item(380, 2).Value = 2024
$ws.Cells.Item(380, 3).Value = 'CDU'
$ws.Cells.Item(380, 4).Value = 'TRUE'
$ws.Cells.Item(380, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/cdusachsenltw-2024.pdf'
$ws.Cells.Item(381, 1).Value = 'SN'
$ws.Cells.Item(381, 2).Value = 2024
$ws.Cells.Item(381, 3).Value = 'Linke'
$ws.Cells.Item(381, 4).Value = $true
$ws.Cells.Item(381, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/dielinkesachsenltw-2024.pdf'
$ws.Cells.Item(382, 1).Value = 'SN'
$ws.Cells.Item(382, 2).Value = 2024
$ws.Cells.Item(382, 3).Value = 'SPD'
$ws.Cells.Item(382, 4).Value = $true
$ws.Cells.Item(382, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdsachsenltw-2024.pdf'
$ws.Cells.Item(383, 1).Value = 'SN'
$ws.Cells.Item(383, 2).Value = 2024
$ws.Cells.Item(383, 3).Value = 'AfD'
$ws.Cells.Item(383, 4).Value = $true
$ws.Cells.Item(383, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afdsachsenltw-2024-1.pdf'
$ws.Cells.Item(384, 1).Value = 'SN'
$ws.Cells.Item(384, 2).Value = 2024
$ws.Cells.Item(384, 3).Value = 'Grüne'
$ws.Cells.Item(384, 4).Value = $true
$ws.Cells.Item(384, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/diegruenensachsenltw-2024.pdf'
$ws.Cells.Item(385, 1).Value = 'SN'
$ws.Cells.Item(385, 2).Value = 2024
$ws.Cells.Item(385, 3).Value = 'FDP'
$ws.Cells.Item(385, 4).Value = $true
$ws.Cells.Item(385, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fdpsachsenltw-2024.pdf'
$ws.Cells.Item(386, 1).Value = 'SN'
$ws.Cells.Item(386, 2).Value = 2024
$ws.Cells.Item(386, 3).Value = 'FW-SN'
$ws.Cells.Item(386, 4).Value = $true
$ws.Cells.Item(386, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/freiewahlersachsenltw-2024.pdf'
$ws.Cells.Item(387, 1).Value = 'SN'
$ws.Cells.Item(387, 2).Value = 2024
$ws.Cells.Item(387, 3).Value = 'ÖDP'
$ws.Cells.Item(387, 4).Value = $true
$ws.Cells.Item(387, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/odpsachsenltw-2024-1.pdf'
$ws.Cells.Item(388, 1).Value = 'SN'
$ws.Cells.Item(388, 2).Value = 2024
$ws.Cells.Item(388, 3).Value = 'BSW'
$ws.Cells.Item(388, 4).Value = $true
$ws.Cells.Item(388, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bswsachsenltw-2024.pdf'
$ws.Cells.Item(389, 1).Value = 'SN'
$ws.Cells.Item(389, 2).Value = 2024
$ws.Cells.Item(389, 3).Value = 'BündnisC'
$ws.Cells.Item(389, 4).Value = $true
$ws.Cells.Item(389, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bundnis-csachsenltw-2024_0.pdf'
$ws.Cells.Item(390, 1).Value = 'SN'
$ws.Cells.Item(390, 2).Value = 2024
$ws.Cells.Item(390, 3).Value = 'BD'
$ws.Cells.Item(390, 4).Value = $true
$ws.Cells.Item(390, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bundnis-deutschlandsachsenltw-2024.pdf'
$ws.Cells.Item(391, 1).Value = 'SN'
$ws.Cells.Item(391, 2).Value = 2024
$ws.Cells.Item(391, 3).Value = 'BÜSO'
$ws.Cells.Item(391, 4).Value = $true
$ws.Cells.Item(391, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/busosachsenltw-2024.pdf'
$ws.Cells.Item(392, 1).Value = 'SN'
$ws.Cells.Item(392, 2).Value = 2024
$ws.Cells.Item(392, 3).Value = 'dieBasis'
$ws.Cells.Item(392, 4).Value = $true
$ws.Cells.Item(392, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/diebasissachsenltw-2024.pdf'
$ws.Cells.Item(393, 1).Value = 'SN'
$ws.Cells.Item(393, 2).Value = 2024
$ws.Cells.Item(393, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/freie-sachsensachsenltw-2024.pdf'
$ws.Cells.Item(393, 3).Value = 'FS (2021)'
$ws.Cells.Item(393, 4).Value = $true
$ws.Cells.Item(394, 1).Value = 'SN'
$ws.Cells.Item(394, 2).Value = 2024
$ws.Cells.Item(394, 3).Value = 'VPartei3'
$ws.Cells.Item(394, 4).Value = $false
$ws.Cells.Item(394, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/v-parteisachsenltw-2024-1.pdf'
$ws.Cells.Item(395, 1).Value = 'SN'
$ws.Cells.Item(395, 2).Value = 2024
$ws.Cells.Item(395, 3).Value = 'Piraten'
$ws.Cells.Item(395, 4).Value = $true
$ws.Cells.Item(395, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/piratensachsenltw-2024.pdf'
$ws.Cells.Item(396, 1).Value = 'TH'
$ws.Cells.Item(396, 2).Value = 2024
$ws.Cells.Item(396, 3).Value = 'Linke'
$ws.Cells.Item(396, 4).Value = $true
$ws.Cells.Item(396, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/linkethuringenltw-2024.pdf'
$ws.Cells.Item(397, 1).Value = 'TH'
$ws.Cells.Item(397, 2).Value = 2024
$ws.Cells.Item(397, 3).Value = 'AfD'
$ws.Cells.Item(397, 4).Value = $true
$ws.Cells.Item(397, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afdthuringenltw-2024.pdf'
$ws.Cells.Item(398, 1).Value = 'TH'
$ws.Cells.Item(398, 2).Value = 2024
$ws.Cells.Item(398, 3).Value = 'CDU'
$ws.Cells.Item(398, 4).Value = $true
$ws.Cells.Item(398, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/cduthuringenltw-2024.pdf'
$ws.Cells.Item(399, 1).Value = 'TH'
$ws.Cells.Item(399, 2).Value = 2024
$ws.Cells.Item(399, 3).Value = 'SPD'
$ws.Cells.Item(399, 4).Value = $true
$ws.Cells.Item(399, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdthuringenltw-2024.pdf'
$ws.Cells.Item(400, 1).Value = 'TH'
$ws.Cells.Item(400, 2).Value = 2024
$ws.Cells.Item(400, 3).Value = 'Grüne'
$ws.Cells.Item(400, 4).Value = $true
$ws.Cells.Item(400, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/die-grunenthuringenltw2024.pdf'
$ws.Cells.Item(401, 1).Value = 'TH'
$ws.Cells.Item(401, 2).Value = 2024
$ws.Cells.Item(401, 3).Value = 'FDP'
$ws.Cells.Item(401, 4).Value = $true
$ws.Cells.Item(401, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fdpthuringenltw-2024.pdf'
$ws.Cells.Item(402, 1).Value = 'TH'
$ws.Cells.Item(402, 2).Value = 2024
$ws.Cells.Item(402, 3).Value = 'TIERSCHUTZliste'
$ws.Cells.Item(402, 4).Value = $true
$ws.Cells.Item(402, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/tierschutz-hierthuringenltw-2024.pdf'
$ws.Cells.Item(403, 1).Value = 'TH'
$ws.Cells.Item(403, 2).Value = 2024
$ws.Cells.Item(403, 3).Value = 'ÖDP'
$ws.Cells.Item(403, 4).Value = $true
$ws.Cells.Item(403, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/odpthuringenltw-2024.pdf'
$ws.Cells.Item(404, 1).Value = 'TH'
$ws.Cells.Item(404, 2).Value = 2024
$ws.Cells.Item(404, 3).Value = 'Piraten'
$ws.Cells.Item(404, 4).Value = $true
$ws.Cells.Item(404, 5).Value = 'https://piraten-thueringen.de/wahlen/wahlprogramm-unsere-leitlinien/'
$ws.Cells.Item(405, 1).Value = 'TH'
$ws.Cells.Item(405, 2).Value = 2024
$ws.Cells.Item(405, 3).Value = 'BSW'
$ws.Cells.Item(405, 4).Value = $true
$ws.Cells.Item(405, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bswthuringenltw-2024.pdf'
$ws.Cells.Item(406, 1).Value = 'TH'
$ws.Cells.Item(406, 2).Value = 2024
$ws.Cells.Item(406, 3).Value = 'FW-TH'
$ws.Cells.Item(406, 4).Value = $true
$ws.Cells.Item(406, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/freie-wahlerthuringenltw-2024.pdf'
$ws.Cells.Item(408, 3).Value = 'WU'
$ws.Cells.Item(407, 1).Value = 'TH'
$ws.Cells.Item(407, 2).Value = 2024
$ws.Cells.Item(407, 3).Value = 'MLPD'
$ws.Cells.Item(407, 4).Value = $true
$ws.Cells.Item(407, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/mlpdthuringenltw-2024.pdf'
$ws.Cells.Item(408, 1).Value = 'TH'
$ws.Cells.Item(408, 2).Value = 2024
$ws.Cells.Item(408, 4).Value = $true
$ws.Cells.Item(408, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/240627-eckpunkte-programm-wu-th-final.pdf'
$ws.Cells.Item(409, 1).Value = 'BB'
$ws.Cells.Item(409, 2).Value = 2024
$ws.Cells.Item(409, 3).Value = 'SPD'
$ws.Cells.Item(409, 4).Value = $true
$ws.Cells.Item(409, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdbrandenburgltw-2024.pdf'
$ws.Cells.Item(410, 1).Value = 'BB'
$ws.Cells.Item(410, 2).Value = 2024
$ws.Cells.Item(410, 3).Value = 'AfD'
$ws.Cells.Item(410, 4).Value = $true
$ws.Cells.Item(410, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afdbrandenburgltw-2024.pdf'
$ws.Cells.Item(411, 1).Value = 'BB'
$ws.Cells.Item(411, 2).Value = 2024
$ws.Cells.Item(411, 3).Value = 'CDU'
$ws.Cells.Item(411, 4).Value = $true
$ws.Cells.Item(411, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/cdubrandenburgltw-2024.pdf'
$ws.Cells.Item(412, 1).Value = 'BB'
$ws.Cells.Item(412, 2).Value = 2024
$ws.Cells.Item(412, 3).Value = 'Grüne'
$ws.Cells.Item(412, 4).Value = $true
$ws.Cells.Item(412, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/grunebrandenburgltw-2024.pdf'
$ws.Cells.Item(413, 1).Value = 'BB'
$ws.Cells.Item(413, 2).Value = 2024
$ws.Cells.Item(413, 3).Value = 'Linke'
$ws.Cells.Item(413, 4).Value = $true
$ws.Cells.Item(413, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/die-linkebrandenburgltw-2024.pdf'
$ws.Cells.Item(414, 1).Value = 'BB'
$ws.Cells.Item(414, 2).Value = 2024
$ws.Cells.Item(414, 3).Value = 'FW-BB'
$ws.Cells.Item(414, 4).Value = $true
$ws.Cells.Item(414, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bvbfwbrandenburgltw-2024.pdf'
$ws.Cells.Item(415, 1).Value = 'BB'
$ws.Cells.Item(415, 2).Value = 2024
$ws.Cells.Item(415, 3).Value = 'FDP'
$ws.Cells.Item(415, 4).Value = $true
$ws.Cells.Item(415, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/2024-07-15-fdp-wahlprogramm-brandenburg-2024.pdf'
$ws.Cells.Item(416, 1).Value = 'BB'
$ws.Cells.Item(416, 2).Value = 2024
$ws.Cells.Item(416, 3).Value = 'DKP'
$ws.Cells.Item(416, 4).Value = $true
$ws.Cells.Item(416, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/dkpbrandenburgltw2024.pdf'
$ws.Cells.Item(417, 1).Value = 'BB'
$ws.Cells.Item(417, 2).Value = 2024
$ws.Cells.Item(417, 3).Value = 'BSW'
$ws.Cells.Item(417, 4).Value = $true
$ws.Cells.Item(417, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bswbrandenburgltw-2024.pdf'
$ws.Cells.Item(418, 1).Value = 'BB'
$ws.Cells.Item(418, 2).Value = 2024
$ws.Cells.Item(418, 3).Value = 'DLW'
$ws.Cells.Item(418, 4).Value = $true
$ws.Cells.Item(418, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/parteiprogramm130624.pdf'
$ws.Cells.Item(419, 1).Value = 'BB'
$ws.Cells.Item(419, 2).Value = 2024
$ws.Cells.Item(419, 3).Value = 'WU'
$ws.Cells.Item(419, 4).Value = $true
$ws.Cells.Item(419, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/unsere-kernthemen-fuer-brandenburg.pdf'
$ws.Cells.Item(420, 1).Value = 'HH'
$ws.Cells.Item(420, 2).Value = 2025
$ws.Cells.Item(420, 3).Value = 'SPD'
$ws.Cells.Item(420, 4).Value = $true
$ws.Cells.Item(420, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdwahlprogramm-buergerschaftswahl-hh-2025pdf.pdf'
$ws.Cells.Item(421, 1).Value = 'HH'
$ws.Cells.Item(421, 2).Value = 2025
$ws.Cells.Item(421, 3).Value = 'Grüne'
$ws.Cells.Item(421, 4).Value = $true
$ws.Cells.Item(421, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/bundnis-90die-grunenwahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(422, 1).Value = 'HH'
$ws.Cells.Item(422, 2).Value = 2025
$ws.Cells.Item(422, 3).Value = 'CDU'
$ws.Cells.Item(422, 4).Value = $true
$ws.Cells.Item(422, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/cduwahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(423, 1).Value = 'HH'
$ws.Cells.Item(423, 2).Value = 2025
$ws.Cells.Item(423, 3).Value = 'Linke'
$ws.Cells.Item(423, 4).Value = $true
$ws.Cells.Item(423, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/die-linkewahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(424, 1).Value = 'HH'
$ws.Cells.Item(424, 2).Value = 2025
$ws.Cells.Item(424, 3).Value = 'AfD'
$ws.Cells.Item(424, 4).Value = $true
$ws.Cells.Item(424, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afdwahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(425, 1).Value = 'HH'
$ws.Cells.Item(425, 2).Value = 2025
$ws.Cells.Item(425, 3).Value = 'FDP'
$ws.Cells.Item(425, 4).Value = $true
$ws.Cells.Item(425, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fdpwahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(426, 1).Value = 'HH'
$ws.Cells.Item(426, 2).Value = 2025
$ws.Cells.Item(426, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/diewahlwahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(426, 3).Value = 'WFG'
$ws.Cells.Item(426, 4).Value = $true
$ws.Cells.Item(427, 1).Value = 'HH'
$ws.Cells.Item(427, 2).Value = 2025
$ws.Cells.Item(427, 3).Value = 'FW-HH'
$ws.Cells.Item(427, 4).Value = $true
$ws.Cells.Item(427, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fwwahlprogramm-burgerschaftswahl-hh-2025.pdf'
$ws.Cells.Item(428, 1).Value = 'HH'
$ws.Cells.Item(428, 2).Value = 2025
$ws.Cells.Item(428, 3).Value = 'Volt'
$ws.Cells.Item(428, 4).Value = $true
$ws.Cells.Item(428, 5).Value = 'https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/voltwahlprogramm-burgerschaftswahl-hh-2025.pdf'

# Row 380 keeps the legacy "TRUE" text convention used by earlier rows in
# this sheet, but picked up an incidental wrap-text toggle (empty <alignment/>).
$ws.Cells.Item(380, 4).WrapText = $false

# Rows 381-428 use real boolean TRUE/FALSE values with left alignment.
for ($r = 381; $r -le 428; $r++) {
    $ws.Cells.Item($r, 4).HorizontalAlignment = -4131
}

# Restore cursor/selection to match the last-edited cell.
$ws.Range("E438").Select()
